$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.711.55'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.26%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.878.21'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.90%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.27%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.10'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.70%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.15%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4736'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.76%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2813'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.93%  '

# Row 9
$ws.Range('E9').Value = '  +3.40%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.50'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +13.68%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.881.32'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.16%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07562'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.96%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '95.15'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +13.61%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.060'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.51%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6473'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.15%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '299.28'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +31.22%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.687.89'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.37%  '

# Row 18
$ws.Range('B18').Value = 'Dai'
$ws.Range('C18').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.001'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.28%  '

# Row 19
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.01'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.37%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007513'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.88%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.130.62'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.93%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('D22').Style = 'Normal'

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.123'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.14%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.128'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.55%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '168.82'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.68%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.184'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.02%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.57'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +9.94%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.937'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.59%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1056'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.43%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.347'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.69%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.140'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.39%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.932'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.30%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05035'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.99%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.165'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.01%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7170'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.02%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.711'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.52%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01908'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.49%  '

# Row 38
$ws.Range('E38').Value = '  +1.71%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.041'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.86%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8941'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.00%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '106.94'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.68%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.001'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.15%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4166'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.61%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.565'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.13%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.287'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.47%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '64.69'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.93%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.936'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.23%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1213'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.33%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '34.46'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.81%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05586'
$ws.Range('D50').Style = 'Normal'

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.369'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.53%  '
